# Insert a new data row at row 21 (shifting existing rows 21..97 down to 22..98)
# and populate it with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(21).Insert()

$ws.Range("A21").Value = 11
$ws.Range("B21").Value = "Vega Monumental Concepción"
$ws.Range("C21").Value = "Bíobío"
$ws.Range("D21").Value2 = 44883
$ws.Range("E21").Value = 8
$ws.Range("F21").Value = "Fruta"
$ws.Range("G21").Value = 100101
$ws.Range("H21").Value = "Berries"
$ws.Range("I21").Value = 100101001
$ws.Range("J21").Value = "Arándano (blue)"
$ws.Range("K21").Value = "Sin especificar"
$ws.Range("L21").Value = "Primera"
$ws.Range("M21").Value = 180
$ws.Range("N21").Value = 6000
$ws.Range("O21").Value = 6500
$ws.Range("P21").Value = 6222
$ws.Range("Q21").Value = "`$/bandeja 2 kilos"
$ws.Range("R21").Value = "Región de O'Higgins"
$ws.Range("S21").Value = 3111
$ws.Range("T21").Value = 2
